$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sayfa1")
$ws2 = $wb.Worksheets.Item("Sayfa2")
$ws3 = $wb.Worksheets.Item("Sayfa3")

# ---------------------------------------------------------------------------
# Sheet1 (Sayfa1): fill in the student info fields
# ---------------------------------------------------------------------------
$ws1.Range("B1").Value = 20215070019
$ws1.Range("B2").Value = "KÜBRA ÇABUK"
$ws1.Range("B3").Value = "YBS"

$ws1.Range("B3:D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet2 (Sayfa2): paste the practice table (letters a..k under column
# headers 2 and 4), re-using the formatting already defined on Sheet1.
# ---------------------------------------------------------------------------
$ws1.Range("J5:K5").Copy()
$ws2.Range("C3:D3").PasteSpecial(-4122)
$ws2.Range("C3").Value = 2
$ws2.Range("D3").Value = 4

$letters = @("a","b","c","d","e","f","g","h","ı","j","k")
for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = 4 + $i
    $ws1.Range("I6").Copy()
    $ws2.Range("B$row").PasteSpecial(-4122)
    $ws2.Range("B$row").Value = $letters[$i]
    $ws1.Range("J6:K6").Copy()
    $ws2.Range("C" + $row + ":D" + $row).PasteSpecial(-4122)
    $ws2.Rows.Item($row).RowHeight = 30
}
$ws2.Rows.Item(3).RowHeight = 30

$ws2.Range("B3:D14").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet3 (Sayfa3): same table after deleting rows b, f, h, j.
# ---------------------------------------------------------------------------
$ws1.Range("J5:K5").Copy()
$ws3.Range("C2:D2").PasteSpecial(-4122)
$ws3.Range("C2").Value = 2
$ws3.Range("D2").Value = 4

$letters3 = @("a","c","d","e","g","ı","k")
for ($i = 0; $i -lt $letters3.Length; $i++) {
    $row = 3 + $i
    $ws1.Range("I6").Copy()
    $ws3.Range("B$row").PasteSpecial(-4122)
    $ws3.Range("B$row").Value = $letters3[$i]
    $ws1.Range("J6:K6").Copy()
    $ws3.Range("C" + $row + ":D" + $row).PasteSpecial(-4122)
    $ws3.Rows.Item($row).RowHeight = 30
}
$ws3.Rows.Item(2).RowHeight = 30

$ws3.Range("G6").Select() | Out-Null

$ws1.Select() | Out-Null
